# Edit script: insert "2022-Q3" quarter data into 06618-JD Health workbook.
# 1. Duplicate the "2022-Q2" sheet (keeps header/column-A styling identical to its
#    siblings), place the copy right after "总计", rename it to "2022-Q3", and overwrite
#    its data with the new quarter's numbers.
# 2. On the "总计" (summary) sheet, shift the existing 7 data rows down by one and write
#    the new 2022-Q3 summary row on top.

$wb = $excel.ActiveWorkbook

# ---- 1. New "2022-Q3" fund-holdings sheet -------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$srcSheet.Copy($null, $totalSheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$raw = @"
164906`t交银施罗德中证海外中国互联网指数（QDII-LOF）`t97.68`t91.19`t3.89`t3.7998`t10
513060`t博时恒生医疗保健ETF（QDII）`t44.00`t99.48`t7.15`t3.1460`t3
159792`t富国中证港股通互联网ETF`t18.08`t99.26`t4.75`t0.8588`t5
517050`t华泰柏瑞中证沪港深互联网ETF`t5.22`t96.35`t3.87`t0.2020`t6
513700`t鹏华中证港股通医药卫生综合ETF`t2.96`t94.14`t6.13`t0.1814`t4
513770`t华宝中证港股通互联网ETF`t3.80`t98.21`t4.68`t0.1778`t5
159856`t工银瑞信中证沪港深互联网ETF`t3.96`t98.04`t3.82`t0.1513`t8
159892`t华夏恒生香港上市生物科技ETF（QDII）`t1.58`t99.13`t7.85`t0.1240`t3
513280`t汇添富恒生香港上市生物科技ETF（QDII）`t1.51`t100.14`t7.95`t0.1200`t3
513860`t海富通中证港股通科技ETF`t3.67`t95.72`t2.75`t0.1009`t10
012379`t创金合信港股互联网3个月持有期混合（QDII）A`t2.81`t87.48`t3.57`t0.1003`t7
517200`t嘉实中证沪港深互联网ETF`t1.42`t97.71`t3.84`t0.0545`t7
513150`t华泰柏瑞中证港股通科技ETF`t1.88`t96.19`t2.70`t0.0508`t10
513200`t易方达中证港股通医药卫生综合ETF`t0.77`t95.67`t6.42`t0.0494`t4
513020`t国泰中证港股通科技ETF`t1.81`t94.40`t2.63`t0.0476`t10
159729`t汇添富中证沪港深互联网ETF`t1.16`t96.75`t3.82`t0.0443`t7
007151`t前海开源沪港深聚瑞混合`t0.53`t82.80`t8.26`t0.0438`t2
003993`t前海开源沪港深核心驱动灵活配置混合`t0.53`t82.41`t7.81`t0.0414`t3
006537`t恒生前海港股通精选混合`t0.95`t90.50`t3.90`t0.0370`t8
012380`t创金合信港股互联网3个月持有期混合（QDII）C`t0.96`t87.48`t3.57`t0.0343`t7
159776`t银华中证港股通医药卫生综合ETF`t0.52`t92.74`t6.04`t0.0314`t4
159718`t平安中证港股通医药卫生综合ETF`t0.53`t90.14`t5.91`t0.0313`t4
004292`t鹏华沪深港互联网股票`t0.77`t92.71`t3.74`t0.0288`t8
159793`t平安中证沪港深线上消费主题ETF`t0.45`t95.11`t4.61`t0.0207`t5
159751`t鹏华中证港股通科技ETF`t0.77`t90.66`t2.50`t0.0192`t10
517280`t天弘中证沪港深线上消费主题ETF`t0.40`t98.98`t4.79`t0.0192`t5
012371`t西藏东财中证沪港深互联网指数A`t0.51`t94.76`t3.72`t0.0190`t7
012372`t西藏东财中证沪港深互联网指数C`t0.44`t94.76`t3.72`t0.0164`t7
006477`t中邮沪港深精选混合`t0.06`t90.21`t4.57`t0.0027`t10
"@
$rows = $raw -split "`n"

# The copied sheet only had 26 data rows (2022-Q2's count); the new quarter needs 29
# (rows 2..30). Extend column A's index styling down to row 30 first (copy forward from a
# row that already carries it) so every index cell stays consistently formatted, matching
# every other quarter sheet in the workbook.
$q3.Range("A26").Copy($q3.Range("A27"))
$q3.Range("A26").Copy($q3.Range("A28"))
$q3.Range("A26").Copy($q3.Range("A29"))
$q3.Range("A26").Copy($q3.Range("A30"))

# Fund code (B) and the scale/position/ratio columns (D:G) are stored as plain text in the
# source data (e.g. "012379" with a leading zero, "44.00"/"3.1460" with trailing zeros) -
# force text format first so COM doesn't silently coerce them to numbers and drop that.
$q3.Range("B2:B30").NumberFormat = "@"
$q3.Range("D2:G30").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $cols = $rows[$i] -split "`t"
    $r = $i + 2
    $q3.Cells.Item($r, 1).Value = $i
    $q3.Cells.Item($r, 2).Value = $cols[0]
    $q3.Cells.Item($r, 3).Value = $cols[1]
    $q3.Cells.Item($r, 4).Value = $cols[2]
    $q3.Cells.Item($r, 5).Value = $cols[3]
    $q3.Cells.Item($r, 6).Value = $cols[4]
    $q3.Cells.Item($r, 7).Value = $cols[5]
    $q3.Cells.Item($r, 8).Value = [double]$cols[6]
}

# ---- 2. "总计" summary sheet: insert the 2022-Q3 row at the top ---------------------------
$ws = $totalSheet

# Give the brand-new row 9 the same column-A styling as the row above it before the shift.
$ws.Range("A8").Copy($ws.Range("A9"))

for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Column A is just the 0-based row index, independent of the shifted quarter data.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

$ws.Cells.Item(2, 2).Value = "2022-Q3"
$ws.Cells.Item(2, 3).Value = 29
$ws.Cells.Item(2, 4).Value = 9.550000000000001
